$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Optimistic")

$row2 = New-Object 'object[,]' 1,28
$row2[0,0] = -8.812779424184445
$row2[0,1] = 17.11756676259009
$row2[0,2] = -7.058713283447916
$row2[0,3] = -8.557142356912657
$row2[0,4] = -3.594502164640295
$row2[0,5] = 10.99744732886962
$row2[0,6] = -16.61037207634313
$row2[0,7] = 5.218337347991672
$row2[0,8] = 0.525139647150322
$row2[0,9] = -6.032539788860974
$row2[0,10] = -25.47886495471957
$row2[0,11] = 8.738558707719651
$row2[0,12] = 3.302091071496324
$row2[0,13] = 9.605454348211776
$row2[0,14] = -9.259411444262277
$row2[0,15] = 19.7699877127707
$row2[0,16] = -9.129912367950297
$row2[0,17] = -13.22284843517242
$row2[0,18] = -6.655821250537621
$row2[0,19] = -13.27355879584231
$row2[0,20] = -14.22625658999419
$row2[0,21] = -6.644221128353055
$row2[0,22] = -5.890039758440476
$row2[0,23] = -14.0566877390623
$row2[0,24] = -35.32259626621381
$row2[0,25] = -21.83004257515117
$row2[0,26] = -17.41989721487705
$row2[0,27] = -47.52835376535336
$ws.Range("B2:AC2").Value = $row2

$row3 = New-Object 'object[,]' 1,28
$row3[0,0] = -5.741461213172084
$row3[0,1] = 1.250483475167204
$row3[0,2] = -13.52450348750358
$row3[0,3] = -13.18389582763616
$row3[0,4] = -7.461952877128537
$row3[0,5] = -12.78440265303405
$row3[0,6] = -5.966039781277625
$row3[0,7] = -8.996582694776825
$row3[0,8] = -14.33726609939034
$row3[0,9] = -18.95871153275746
$row3[0,10] = -14.42501588398642
$row3[0,11] = -6.852405415713553
$row3[0,12] = -13.02532793824343
$row3[0,13] = 3.677201255883681
$row3[0,14] = -19.32398865577915
$row3[0,15] = -27.93768776871568
$row3[0,16] = -3.777657297778002
$row3[0,17] = -28.98481682210415
$row3[0,18] = -19.14776711416241
$row3[0,19] = -33.05396947232694
$row3[0,20] = -16.59139585123468
$row3[0,21] = -21.72572513394365
$row3[0,22] = -21.04043228443193
$row3[0,23] = -23.54569298786582
$row3[0,24] = -22.98266608733261
$row3[0,25] = -22.0146684999642
$row3[0,26] = -42.76753636561547
$row3[0,27] = -47.75328740501817
$ws.Range("B3:AC3").Value = $row3

$row4 = New-Object 'object[,]' 1,28
$row4[0,0] = 6.022772527081343
$row4[0,1] = 6.976168701760161
$row4[0,2] = -2.447248192114507
$row4[0,3] = 4.408413611006896
$row4[0,4] = -19.56551522344432
$row4[0,5] = -8.502105137517841
$row4[0,6] = -11.72609121189239
$row4[0,7] = -10.33343157800917
$row4[0,8] = -0.5609008847496169
$row4[0,9] = -7.804422357898893
$row4[0,10] = -6.764455079141137
$row4[0,11] = -21.01380300536228
$row4[0,12] = -3.672479346840388
$row4[0,13] = -9.293895538956981
$row4[0,14] = -1.143465379249979
$row4[0,15] = 0.3529789371306578
$row4[0,16] = -10.42292625476901
$row4[0,17] = 2.358552203403535
$row4[0,18] = -19.75747571287381
$row4[0,19] = -23.15147422497265
$row4[0,20] = -10.49223775469964
$row4[0,21] = -11.88585340621532
$row4[0,22] = -13.2936323404691
$row4[0,23] = -31.16314101429282
$row4[0,24] = -31.47363410395297
$row4[0,25] = -6.625514569927279
$row4[0,26] = -52.92811036104301
$row4[0,27] = -46.85255844873348
$ws.Range("B4:AC4").Value = $row4

$row5 = New-Object 'object[,]' 1,28
$row5[0,0] = 2.096962369033812
$row5[0,1] = 14.57555102440856
$row5[0,2] = -23.5745591773492
$row5[0,3] = -1.033559171613133
$row5[0,4] = -15.91758556250333
$row5[0,5] = -2.57199913612783
$row5[0,6] = -16.04550684505481
$row5[0,7] = -5.96923248036609
$row5[0,8] = -9.110304262026849
$row5[0,9] = -23.39711590823227
$row5[0,10] = -5.479627662700837
$row5[0,11] = -8.235865476313478
$row5[0,12] = -8.948366307192604
$row5[0,13] = -21.67989243069401
$row5[0,14] = -21.49112010601729
$row5[0,15] = -4.176915716718425
$row5[0,16] = 2.133818501362063
$row5[0,17] = -5.426947188468064
$row5[0,18] = -8.392782186380437
$row5[0,19] = -13.16378559554878
$row5[0,20] = -14.37318454791788
$row5[0,21] = -20.0509995774255
$row5[0,22] = -32.2639725421937
$row5[0,23] = -21.37721676942371
$row5[0,24] = -16.62468406663643
$row5[0,25] = -34.65360233279628
$row5[0,26] = -23.91687352070372
$row5[0,27] = -36.79521894131058
$ws.Range("B5:AC5").Value = $row5

$row6 = New-Object 'object[,]' 1,28
$row6[0,0] = 7.135814955431026
$row6[0,1] = 5.56044112588655
$row6[0,2] = -24.30612410185532
$row6[0,3] = -6.313216051853364
$row6[0,4] = 6.472574054261148
$row6[0,5] = 2.929975878031675
$row6[0,6] = -18.95797203404481
$row6[0,7] = -9.908327420051155
$row6[0,8] = 24.19563821344937
$row6[0,9] = -26.02303599435923
$row6[0,10] = -18.91498097573676
$row6[0,11] = -3.633567507636424
$row6[0,12] = 7.901911128187672
$row6[0,13] = -14.03003306103507
$row6[0,14] = -0.2562100030331109
$row6[0,15] = -13.11333180700851
$row6[0,16] = -20.59626496618199
$row6[0,17] = -18.26400392199824
$row6[0,18] = -10.89527507925232
$row6[0,19] = 3.451104592111074
$row6[0,20] = -14.59444872892221
$row6[0,21] = -22.75047699907367
$row6[0,22] = -15.78606036267524
$row6[0,23] = -20.08673377348418
$row6[0,24] = -7.565154646360496
$row6[0,25] = 0.9378814045714599
$row6[0,26] = -28.85101561615434
$row6[0,27] = -52.62002631754913
$ws.Range("B6:AC6").Value = $row6

$row7 = New-Object 'object[,]' 1,28
$row7[0,0] = 8.796666149667619
$row7[0,1] = 5.268555221427171
$row7[0,2] = -19.14724627072031
$row7[0,3] = 0.8374461163785805
$row7[0,4] = 0.6609881448083961
$row7[0,5] = -4.780525359676183
$row7[0,6] = -5.063827041921656
$row7[0,7] = -16.99681034156701
$row7[0,8] = 6.355972393606649
$row7[0,9] = -15.21853197419455
$row7[0,10] = -4.794189392778666
$row7[0,11] = -3.783132834497415
$row7[0,12] = -3.0349358142398
$row7[0,13] = -0.4981420101134679
$row7[0,14] = 7.01864943365276
$row7[0,15] = 1.189513772469379
$row7[0,16] = -19.40599719847911
$row7[0,17] = -4.61587994828213
$row7[0,18] = -16.23149027048957
$row7[0,19] = -1.315656117354049
$row7[0,20] = -23.74967368708842
$row7[0,21] = -4.859122061115725
$row7[0,22] = -16.97843307202257
$row7[0,23] = -39.56701082707624
$row7[0,24] = -23.26056358901563
$row7[0,25] = -5.252824718505618
$row7[0,26] = -28.54647792440355
$row7[0,27] = -46.10307480874474
$ws.Range("B7:AC7").Value = $row7

$row8 = New-Object 'object[,]' 1,28
$row8[0,0] = -6.503055484048042
$row8[0,1] = 3.410204092652535
$row8[0,2] = -13.60521653594496
$row8[0,3] = -4.865481326417701
$row8[0,4] = -21.45838570904424
$row8[0,5] = -13.39909843819251
$row8[0,6] = -20.68034006856332
$row8[0,7] = -23.71172804470534
$row8[0,8] = 22.89649497897013
$row8[0,9] = -9.865883255861675
$row8[0,10] = -1.860075801079909
$row8[0,11] = 1.339787468278868
$row8[0,12] = -15.56988817448469
$row8[0,13] = -6.165700901007051
$row8[0,14] = -33.37841804739845
$row8[0,15] = 5.160379188073643
$row8[0,16] = -22.79611594594778
$row8[0,17] = -17.4536107322315
$row8[0,18] = -16.85978909655059
$row8[0,19] = -20.43972447335449
$row8[0,20] = -1.154100944017992
$row8[0,21] = -31.50353870547752
$row8[0,22] = -26.73711804179469
$row8[0,23] = -23.77685665758916
$row8[0,24] = -23.50642080077673
$row8[0,25] = -3.925871058089328
$row8[0,26] = -30.21473839097254
$row8[0,27] = -60.1822604868297
$ws.Range("B8:AC8").Value = $row8

$row9 = New-Object 'object[,]' 1,28
$row9[0,0] = -1.679403765209468
$row9[0,1] = 15.24683649981244
$row9[0,2] = -0.2160900397866339
$row9[0,3] = 2.77963089150523
$row9[0,4] = 9.572928907133624
$row9[0,5] = -17.79339323225614
$row9[0,6] = -1.680137652675097
$row9[0,7] = -1.924447945238266
$row9[0,8] = -26.47461480123049
$row9[0,9] = -4.705903887903705
$row9[0,10] = -9.555874570218897
$row9[0,11] = -12.72654122833984
$row9[0,12] = 0.5273448630426198
$row9[0,13] = -12.75350153142142
$row9[0,14] = -4.261498304332386
$row9[0,15] = 4.567488117211578
$row9[0,16] = 4.462154535301048
$row9[0,17] = -5.200443228670458
$row9[0,18] = -0.6590800824816476
$row9[0,19] = -12.41800917122425
$row9[0,20] = -8.19097373794386
$row9[0,21] = -22.33377106455513
$row9[0,22] = -15.32646191027227
$row9[0,23] = -30.60768759773033
$row9[0,24] = -20.46411527554369
$row9[0,25] = -23.22764307299513
$row9[0,26] = -31.97048834550982
$row9[0,27] = -41.01208512287833
$ws.Range("B9:AC9").Value = $row9

$row10 = New-Object 'object[,]' 1,28
$row10[0,0] = -19.79819545194769
$row10[0,1] = -6.902620810424926
$row10[0,2] = -4.76221881550526
$row10[0,3] = 0.4223071288471765
$row10[0,4] = 2.849102781583942
$row10[0,5] = 2.28912769641107
$row10[0,6] = -25.73666185244208
$row10[0,7] = -5.330575029396122
$row10[0,8] = -21.22193682019407
$row10[0,9] = -6.303878829915901
$row10[0,10] = -20.90383080263397
$row10[0,11] = -19.63985550215957
$row10[0,12] = -14.53743141211644
$row10[0,13] = -30.16439051557035
$row10[0,14] = -24.08480920057367
$row10[0,15] = -21.10283190301825
$row10[0,16] = -7.272683277192602
$row10[0,17] = -9.698866220380282
$row10[0,18] = -25.07942328693183
$row10[0,19] = -6.065514002900872
$row10[0,20] = -16.48770860505826
$row10[0,21] = -27.69687227940149
$row10[0,22] = -17.41143834635319
$row10[0,23] = -7.825362645180187
$row10[0,24] = -25.34189778581194
$row10[0,25] = -20.06159077635053
$row10[0,26] = -16.46329472551048
$row10[0,27] = -56.50039570520566
$ws.Range("B10:AC10").Value = $row10

$row11 = New-Object 'object[,]' 1,28
$row11[0,0] = 1.283811050688735
$row11[0,1] = -8.244503331564836
$row11[0,2] = -24.92069436001444
$row11[0,3] = -5.491207259887902
$row11[0,4] = -4.598493789137732
$row11[0,5] = -13.33353959572788
$row11[0,6] = -9.275083549139813
$row11[0,7] = -4.528529035564008
$row11[0,8] = -11.28418063421509
$row11[0,9] = -2.841061189802809
$row11[0,10] = -19.5033934991747
$row11[0,11] = -9.662547941524112
$row11[0,12] = 0.6583618569473071
$row11[0,13] = 0.4949692679927216
$row11[0,14] = 1.586705627760564
$row11[0,15] = -1.465949834745182
$row11[0,16] = -7.956253408544852
$row11[0,17] = -31.53236511050618
$row11[0,18] = -5.369294585568845
$row11[0,19] = -1.083194744395557
$row11[0,20] = -23.33229816093734
$row11[0,21] = -22.05189009004159
$row11[0,22] = -24.23231031288725
$row11[0,23] = -27.09881683208349
$row11[0,24] = -24.47042017577613
$row11[0,25] = -19.90809009207921
$row11[0,26] = -38.58700187047955
$row11[0,27] = -64.18381610591747
$ws.Range("B11:AC11").Value = $row11

$ws = $wb.Worksheets.Item("Pessimistic")

$row2 = New-Object 'object[,]' 1,28
$row2[0,0] = -8.445699459704631
$row2[0,1] = -4.862194063725537
$row2[0,2] = -15.89688220623006
$row2[0,3] = 11.77947903359717
$row2[0,4] = -6.66564250290854
$row2[0,5] = 2.563856657654434
$row2[0,6] = -7.525171312768133
$row2[0,7] = 1.078057486706943
$row2[0,8] = 4.671951277104801
$row2[0,9] = 6.33926806330791
$row2[0,10] = 16.66427753005458
$row2[0,11] = -3.522338165741361
$row2[0,12] = -5.447083475939571
$row2[0,13] = -1.232763537488992
$row2[0,14] = -2.259433602552217
$row2[0,15] = 19.39683600957655
$row2[0,16] = 9.205135260614771
$row2[0,17] = -6.14441831949676
$row2[0,18] = 11.18054610266636
$row2[0,19] = 1.651402166182923
$row2[0,20] = -8.056800272561631
$row2[0,21] = -0.1033516370600798
$row2[0,22] = -13.85436283191828
$row2[0,23] = 3.492158701691356
$row2[0,24] = -0.4898191886065604
$row2[0,25] = -4.144477416462786
$row2[0,26] = -8.076866511191284
$row2[0,27] = 1.967220948076386
$ws.Range("B2:AC2").Value = $row2

$row3 = New-Object 'object[,]' 1,28
$row3[0,0] = -4.703918524380742
$row3[0,1] = -7.384701731819304
$row3[0,2] = -17.12908476817802
$row3[0,3] = 13.1006902674003
$row3[0,4] = -8.919978352751468
$row3[0,5] = -7.793667272376371
$row3[0,6] = -15.46447442159888
$row3[0,7] = -7.596156118479681
$row3[0,8] = -16.6511084491566
$row3[0,9] = 2.329175800029099
$row3[0,10] = -15.42391788323441
$row3[0,11] = -1.772978384570476
$row3[0,12] = -2.279131274530956
$row3[0,13] = -11.49730449259545
$row3[0,14] = 9.172175386560804
$row3[0,15] = 11.15910059004865
$row3[0,16] = 0.3791650468645855
$row3[0,17] = -5.518081377613904
$row3[0,18] = 1.194147547211803
$row3[0,19] = -3.464616300116811
$row3[0,20] = -11.71753612733431
$row3[0,21] = 7.210459328607246
$row3[0,22] = -7.309279403474726
$row3[0,23] = -1.74135169750864
$row3[0,24] = 6.953424874641135
$row3[0,25] = -9.566716004939725
$row3[0,26] = -9.829166651699715
$row3[0,27] = -3.294978376047523
$ws.Range("B3:AC3").Value = $row3

$row4 = New-Object 'object[,]' 1,28
$row4[0,0] = 1.96788411156886
$row4[0,1] = 9.793360833426242
$row4[0,2] = 2.449572975446845
$row4[0,3] = 19.70452590623813
$row4[0,4] = -18.00356596829654
$row4[0,5] = 1.060555171804593
$row4[0,6] = -8.927318896318944
$row4[0,7] = 0.7766874153775793
$row4[0,8] = -21.40813162777153
$row4[0,9] = -5.750874529090413
$row4[0,10] = -11.29953055288927
$row4[0,11] = -11.74913540106592
$row4[0,12] = -9.687024275935931
$row4[0,13] = 21.44633726923782
$row4[0,14] = -4.333644245300228
$row4[0,15] = 14.15437335288479
$row4[0,16] = -2.477744217123757
$row4[0,17] = 0.7085105858664238
$row4[0,18] = -15.07960216819733
$row4[0,19] = 6.456702166358234
$row4[0,20] = 0.6040383845647335
$row4[0,21] = 8.524903913752267
$row4[0,22] = -11.25129838446836
$row4[0,23] = -6.163730165653814
$row4[0,24] = 18.99136237158999
$row4[0,25] = -6.320691244116416
$row4[0,26] = 7.458478334669625
$row4[0,27] = 1.350437266684461
$ws.Range("B4:AC4").Value = $row4

$row5 = New-Object 'object[,]' 1,28
$row5[0,0] = 5.69540421259098
$row5[0,1] = 15.8443387737497
$row5[0,2] = -9.041383510662909
$row5[0,3] = -1.807895692866627
$row5[0,4] = -22.93028848654954
$row5[0,5] = -3.67392416836885
$row5[0,6] = -19.0499390696505
$row5[0,7] = 2.007711050085545
$row5[0,8] = -20.86281846529436
$row5[0,9] = 2.397143887738131
$row5[0,10] = -8.921630328835224
$row5[0,11] = -12.69261308317471
$row5[0,12] = 7.454930149248398
$row5[0,13] = 13.8524889489197
$row5[0,14] = 4.97633624219995
$row5[0,15] = -3.269501203145849
$row5[0,16] = -2.254427120216117
$row5[0,17] = -22.30869218413394
$row5[0,18] = 5.404724735478939
$row5[0,19] = -14.55394099190735
$row5[0,20] = 4.267262617854231
$row5[0,21] = 8.569253848633311
$row5[0,22] = -11.05058435546473
$row5[0,23] = -7.364367792150111
$row5[0,24] = 7.071670998689303
$row5[0,25] = 14.14281613981135
$row5[0,26] = -11.20883218985977
$row5[0,27] = 0.04303784302254421
$ws.Range("B5:AC5").Value = $row5

$row6 = New-Object 'object[,]' 1,28
$row6[0,0] = 3.179892551929328
$row6[0,1] = 7.304041360956223
$row6[0,2] = 2.139503181013015
$row6[0,3] = 21.32115225012946
$row6[0,4] = -6.922125838414118
$row6[0,5] = -16.90753063634777
$row6[0,6] = 8.67484141237636
$row6[0,7] = -9.511256722667447
$row6[0,8] = -13.27015633037353
$row6[0,9] = -1.125013270616859
$row6[0,10] = -10.67015299958825
$row6[0,11] = -10.58432324433678
$row6[0,12] = -11.15310070730012
$row6[0,13] = -1.360848297030258
$row6[0,14] = -10.31030141346771
$row6[0,15] = 0.4360796688477184
$row6[0,16] = 16.38307383072669
$row6[0,17] = -3.628815292795601
$row6[0,18] = 20.16238805919677
$row6[0,19] = 2.975537069133482
$row6[0,20] = 3.318608173135195
$row6[0,21] = 4.287128448840129
$row6[0,22] = 4.297348869863408
$row6[0,23] = 13.36251912857334
$row6[0,24] = -1.325818190106531
$row6[0,25] = 1.067634855129112
$row6[0,26] = -0.4388569801243412
$row6[0,27] = -5.850561563560658
$ws.Range("B6:AC6").Value = $row6

$row7 = New-Object 'object[,]' 1,28
$row7[0,0] = -11.5491570097279
$row7[0,1] = -3.164142161908076
$row7[0,2] = -0.09112492226020219
$row7[0,3] = -0.4896661936531359
$row7[0,4] = -2.965875791902664
$row7[0,5] = -18.64874809014069
$row7[0,6] = -5.363887656098617
$row7[0,7] = 24.15670224402538
$row7[0,8] = -12.30721989073341
$row7[0,9] = -3.420981975518997
$row7[0,10] = -5.738221092061988
$row7[0,11] = -2.333185794503485
$row7[0,12] = 1.841576508904545
$row7[0,13] = 3.647772968113578
$row7[0,14] = -6.750752307649845
$row7[0,15] = -5.484008345445941
$row7[0,16] = -1.787946318722657
$row7[0,17] = 7.847697987326093
$row7[0,18] = -0.700916487006209
$row7[0,19] = -18.02223587003803
$row7[0,20] = -5.703493722650853
$row7[0,21] = -10.93386618612735
$row7[0,22] = 2.191117502207945
$row7[0,23] = 4.260841497552554
$row7[0,24] = 19.29557750612478
$row7[0,25] = 13.36508815014068
$row7[0,26] = -29.7526954062402
$row7[0,27] = -6.410806346881802
$ws.Range("B7:AC7").Value = $row7

$row8 = New-Object 'object[,]' 1,28
$row8[0,0] = -9.487714368200098
$row8[0,1] = -10.54629348498882
$row8[0,2] = -18.01023692089304
$row8[0,3] = 0.5010206421787267
$row8[0,4] = 3.520693305133928
$row8[0,5] = -16.93469311297584
$row8[0,6] = -4.19833907938612
$row8[0,7] = 2.131827916228263
$row8[0,8] = -1.307299684998619
$row8[0,9] = 15.14158828693672
$row8[0,10] = -3.230046673618826
$row8[0,11] = -25.35596339003044
$row8[0,12] = 12.63027632621623
$row8[0,13] = -23.43030145702912
$row8[0,14] = 10.86481008290672
$row8[0,15] = 4.128037866053131
$row8[0,16] = 3.820237975486838
$row8[0,17] = -16.09448771596702
$row8[0,18] = 9.895866262690834
$row8[0,19] = 2.317334035479305
$row8[0,20] = -1.253336303752502
$row8[0,21] = -15.29957968862895
$row8[0,22] = -3.679666939968412
$row8[0,23] = -4.03339461842495
$row8[0,24] = 6.703046720289267
$row8[0,25] = -9.617241695780034
$row8[0,26] = -18.55429033488829
$row8[0,27] = -9.630250198531394
$ws.Range("B8:AC8").Value = $row8

$row9 = New-Object 'object[,]' 1,28
$row9[0,0] = -1.898984775211348
$row9[0,1] = -5.050551684939466
$row9[0,2] = -9.504131979289451
$row9[0,3] = 10.59053926534656
$row9[0,4] = -18.134787500035
$row9[0,5] = -6.773941278092058
$row9[0,6] = 3.985371424650133
$row9[0,7] = 2.257050013888467
$row9[0,8] = -3.131228081658928
$row9[0,9] = 7.262278955497596
$row9[0,10] = 9.49733374703329
$row9[0,11] = -4.908057935475792
$row9[0,12] = 13.57883974637917
$row9[0,13] = -5.784824944697171
$row9[0,14] = -8.356984217585614
$row9[0,15] = 4.775874178313871
$row9[0,16] = 17.37444926907451
$row9[0,17] = -21.88207193513947
$row9[0,18] = -3.098411724621789
$row9[0,19] = 2.89438987340292
$row9[0,20] = -17.32848166671971
$row9[0,21] = 26.03938190446221
$row9[0,22] = -3.298573616297135
$row9[0,23] = -7.960698470846839
$row9[0,24] = -7.734410007691003
$row9[0,25] = -15.34293803644868
$row9[0,26] = 1.937934573012914
$row9[0,27] = -3.218760903742135
$ws.Range("B9:AC9").Value = $row9

$row10 = New-Object 'object[,]' 1,28
$row10[0,0] = -5.417532623322577
$row10[0,1] = -11.69386591171691
$row10[0,2] = -16.1455796190507
$row10[0,3] = -0.2707313298397902
$row10[0,4] = -1.885623856261962
$row10[0,5] = 3.558987314318966
$row10[0,6] = -1.305763515458318
$row10[0,7] = -0.09936650621417309
$row10[0,8] = -19.67477106923535
$row10[0,9] = 0.2850805195645716
$row10[0,10] = -11.52912951906626
$row10[0,11] = -13.89887660964032
$row10[0,12] = -10.88851491126718
$row10[0,13] = 2.845491638898932
$row10[0,14] = 11.58799355520798
$row10[0,15] = 14.88826541100849
$row10[0,16] = 13.62925126979036
$row10[0,17] = 10.23127818649927
$row10[0,18] = -3.55600514412768
$row10[0,19] = 1.976515574589449
$row10[0,20] = -0.8952782819515051
$row10[0,21] = -10.22231149647097
$row10[0,22] = -13.99694840141624
$row10[0,23] = -9.380937028547269
$row10[0,24] = -9.990080843538729
$row10[0,25] = 2.307297557052507
$row10[0,26] = 1.154772315805278
$row10[0,27] = -3.211795851472504
$ws.Range("B10:AC10").Value = $row10

$row11 = New-Object 'object[,]' 1,28
$row11[0,0] = -5.803780916214143
$row11[0,1] = 2.30433470568707
$row11[0,2] = 5.894405827213402
$row11[0,3] = 6.607366362117334
$row11[0,4] = 0.1795947928747172
$row11[0,5] = -9.946795262272431
$row11[0,6] = -10.84170992544662
$row11[0,7] = 0.09090767135536426
$row11[0,8] = -17.4453570347133
$row11[0,9] = 16.18535308066385
$row11[0,10] = -13.41977381348576
$row11[0,11] = -6.960191808726066
$row11[0,12] = -15.58782468688498
$row11[0,13] = 19.47823103185367
$row11[0,14] = -10.39249478609842
$row11[0,15] = 1.679285639895129
$row11[0,16] = 5.792273192994559
$row11[0,17] = -6.428430438611606
$row11[0,18] = -3.290580108857488
$row11[0,19] = -4.01387247521526
$row11[0,20] = -9.774757416276827
$row11[0,21] = -8.087235273940696
$row11[0,22] = -10.80958925113302
$row11[0,23] = 18.99561005351461
$row11[0,24] = -2.703777173911718
$row11[0,25] = 0.5165469295324172
$row11[0,26] = -2.938647882657365
$row11[0,27] = 24.14875703831298
$ws.Range("B11:AC11").Value = $row11

$ws = $wb.Worksheets.Item("Middle")

$row2 = New-Object 'object[,]' 1,28
$row2[0,0] = -11.82067171103876
$row2[0,1] = 8.726748548575348
$row2[0,2] = 4.267524925961021
$row2[0,3] = -0.6536713277386399
$row2[0,4] = 4.119206919697215
$row2[0,5] = -0.7704658907991684
$row2[0,6] = 2.930229259348289
$row2[0,7] = -19.74581047170735
$row2[0,8] = 7.625541190602304
$row2[0,9] = 7.175884822735918
$row2[0,10] = -20.63894792286762
$row2[0,11] = 3.881828170585614
$row2[0,12] = -15.49147567784482
$row2[0,13] = 8.079012344252993
$row2[0,14] = -13.25388363211106
$row2[0,15] = -14.18745537687296
$row2[0,16] = 7.263392546524384
$row2[0,17] = -15.77921126800205
$row2[0,18] = 4.320778800946329
$row2[0,19] = -11.47399778578884
$row2[0,20] = -10.97064856577721
$row2[0,21] = -1.857933487979964
$row2[0,22] = -1.429640772471988
$row2[0,23] = -5.944332643392298
$row2[0,24] = 4.175096498703073
$row2[0,25] = -28.39176124135981
$row2[0,26] = -8.803471318909994
$row2[0,27] = -41.39606310740389
$ws.Range("B2:AC2").Value = $row2

$row3 = New-Object 'object[,]' 1,28
$row3[0,0] = -16.94790968237312
$row3[0,1] = -11.62188119921013
$row3[0,2] = -5.358798950098826
$row3[0,3] = -4.683155295688005
$row3[0,4] = -14.42090544026057
$row3[0,5] = -15.67730812338714
$row3[0,6] = -17.98741896048194
$row3[0,7] = -7.834353487757377
$row3[0,8] = -4.9669640013528
$row3[0,9] = -28.75006310376186
$row3[0,10] = 2.031988676315429
$row3[0,11] = 15.50820611764182
$row3[0,12] = 0.6109380571451677
$row3[0,13] = -10.69841352279428
$row3[0,14] = -21.86771121364729
$row3[0,15] = -10.3821520385503
$row3[0,16] = 3.4819536462878
$row3[0,17] = -6.601069874086791
$row3[0,18] = -1.648739501752539
$row3[0,19] = -7.932649576409164
$row3[0,20] = -29.00389214601918
$row3[0,21] = -13.75756075645067
$row3[0,22] = -2.585793525559178
$row3[0,23] = -7.676320221983371
$row3[0,24] = -10.08209651484703
$row3[0,25] = -22.16813334412639
$row3[0,26] = -4.600520653722944
$row3[0,27] = -31.07035354616133
$ws.Range("B3:AC3").Value = $row3

$row4 = New-Object 'object[,]' 1,28
$row4[0,0] = -4.839513657942192
$row4[0,1] = 6.53967324899368
$row4[0,2] = 19.11750044879022
$row4[0,3] = -10.31826086516283
$row4[0,4] = -6.806026926225879
$row4[0,5] = -25.0526589369119
$row4[0,6] = -21.21735327493536
$row4[0,7] = 2.911602225575698
$row4[0,8] = 9.620448230639152
$row4[0,9] = 15.97331702273044
$row4[0,10] = -0.400615693395614
$row4[0,11] = -12.25883270681523
$row4[0,12] = -13.11748496073857
$row4[0,13] = -0.1392639755792497
$row4[0,14] = -5.977484828873105
$row4[0,15] = -19.44668051854613
$row4[0,16] = -8.80673940794882
$row4[0,17] = -15.97015419744475
$row4[0,18] = -2.597338392924615
$row4[0,19] = -3.092388827180172
$row4[0,20] = 8.443310182828382
$row4[0,21] = 12.75622526377611
$row4[0,22] = 3.176632091714088
$row4[0,23] = -20.5756990460565
$row4[0,24] = 0.1350901041329546
$row4[0,25] = -5.121551253197074
$row4[0,26] = -11.45179433426688
$row4[0,27] = -35.66936050847141
$ws.Range("B4:AC4").Value = $row4

$row5 = New-Object 'object[,]' 1,28
$row5[0,0] = 0.8564241484497042
$row5[0,1] = -7.476852344089739
$row5[0,2] = -12.49013320649768
$row5[0,3] = 6.107314116929471
$row5[0,4] = -8.903138175978349
$row5[0,5] = -16.91077507026007
$row5[0,6] = -15.47305673262009
$row5[0,7] = -1.515991562119174
$row5[0,8] = -1.427247358790395
$row5[0,9] = -22.76224368886097
$row5[0,10] = 15.06103336748857
$row5[0,11] = -17.52531352008921
$row5[0,12] = -15.4839363322223
$row5[0,13] = -7.994138963824208
$row5[0,14] = -14.12427060111342
$row5[0,15] = -9.974912337127002
$row5[0,16] = 7.831240285440626
$row5[0,17] = -7.589247767652509
$row5[0,18] = -11.3840530590101
$row5[0,19] = -4.507649000703443
$row5[0,20] = -11.87654273985714
$row5[0,21] = -13.58343908620171
$row5[0,22] = 8.27726142643284
$row5[0,23] = -21.36000318293432
$row5[0,24] = 13.11023386972459
$row5[0,25] = -16.74469773150762
$row5[0,26] = -16.90650244946344
$row5[0,27] = -54.68863172876889
$ws.Range("B5:AC5").Value = $row5

$row6 = New-Object 'object[,]' 1,28
$row6[0,0] = 7.651614100858552
$row6[0,1] = 2.982350242261731
$row6[0,2] = -9.423333829321006
$row6[0,3] = -2.123044486343237
$row6[0,4] = 10.52802277054863
$row6[0,5] = -7.064597153999751
$row6[0,6] = 12.70408493941412
$row6[0,7] = -0.771879458547382
$row6[0,8] = 7.55024441805049
$row6[0,9] = -6.732961033964751
$row6[0,10] = 1.01804234719919
$row6[0,11] = 6.049847242443797
$row6[0,12] = -19.79404415843948
$row6[0,13] = 18.19490563739133
$row6[0,14] = 16.65774256994644
$row6[0,15] = -9.963191788770876
$row6[0,16] = -8.739759981814199
$row6[0,17] = 0.7120890284270835
$row6[0,18] = -18.61659693709531
$row6[0,19] = 1.768827057095822
$row6[0,20] = 9.178223178747952
$row6[0,21] = 0.1670004055298744
$row6[0,22] = 11.98721579785152
$row6[0,23] = -0.1143203181843759
$row6[0,24] = -4.017185308309948
$row6[0,25] = -27.4778044891386
$row6[0,26] = 3.514931842842902
$row6[0,27] = -28.30841831843108
$ws.Range("B6:AC6").Value = $row6

$row7 = New-Object 'object[,]' 1,28
$row7[0,0] = 3.435242626192617
$row7[0,1] = -6.211067367464246
$row7[0,2] = 2.428545967542385
$row7[0,3] = 2.845883260940616
$row7[0,4] = 1.974166652740444
$row7[0,5] = -24.05439544074044
$row7[0,6] = -13.9390486926758
$row7[0,7] = 3.938970407756906
$row7[0,8] = -5.411435738014069
$row7[0,9] = 14.03817134799943
$row7[0,10] = -12.23764753788003
$row7[0,11] = -11.83969584341656
$row7[0,12] = 0.7780602366001466
$row7[0,13] = -8.661515140153124
$row7[0,14] = -3.739563792826059
$row7[0,15] = -3.781839073044715
$row7[0,16] = 28.66428693032623
$row7[0,17] = -17.2511790601209
$row7[0,18] = 17.99735417640368
$row7[0,19] = 14.64027564218928
$row7[0,20] = -6.760386987332817
$row7[0,21] = 5.550619223956054
$row7[0,22] = -17.66053870660931
$row7[0,23] = 1.747568439450616
$row7[0,24] = -1.746637291183712
$row7[0,25] = -7.447814022222066
$row7[0,26] = 6.502737954629477
$row7[0,27] = -44.61830314407507
$ws.Range("B7:AC7").Value = $row7

$row8 = New-Object 'object[,]' 1,28
$row8[0,0] = -15.90971247723647
$row8[0,1] = -5.627722207001062
$row8[0,2] = -12.32611370703521
$row8[0,3] = -0.3795358579761849
$row8[0,4] = 8.115925983291984
$row8[0,5] = -2.260290485882043
$row8[0,6] = -0.2571779570006747
$row8[0,7] = -10.3045238777016
$row8[0,8] = -11.69028729091807
$row8[0,9] = -8.196146027751176
$row8[0,10] = 8.770715706558148
$row8[0,11] = -3.814599226789914
$row8[0,12] = -15.32574933173023
$row8[0,13] = 15.14069079071709
$row8[0,14] = 3.841486808113955
$row8[0,15] = -13.6211795266525
$row8[0,16] = -8.340118181411874
$row8[0,17] = 16.88046091869175
$row8[0,18] = -17.84889762987029
$row8[0,19] = -7.90959243971433
$row8[0,20] = -8.300117785792489
$row8[0,21] = 4.333612722154336
$row8[0,22] = -16.77669573276757
$row8[0,23] = -20.11392220029712
$row8[0,24] = -1.187253132140917
$row8[0,25] = -26.06129251246172
$row8[0,26] = -13.00989373662897
$row8[0,27] = -46.79870026326006
$ws.Range("B8:AC8").Value = $row8

$row9 = New-Object 'object[,]' 1,28
$row9[0,0] = -6.398822348157381
$row9[0,1] = -1.632516700766025
$row9[0,2] = 5.376492440797194
$row9[0,3] = 8.296867895445702
$row9[0,4] = 13.19860333908353
$row9[0,5] = -6.196503588947222
$row9[0,6] = -5.045851777839451
$row9[0,7] = -2.326403830537651
$row9[0,8] = 1.469569212263153
$row9[0,9] = -22.9771658577931
$row9[0,10] = 6.352417926464163
$row9[0,11] = 2.200206823091766
$row9[0,12] = -13.41507458488935
$row9[0,13] = -3.870263471076483
$row9[0,14] = 6.076644119401887
$row9[0,15] = -9.920983542701988
$row9[0,16] = 8.142713002255434
$row9[0,17] = 5.247126404065011
$row9[0,18] = -3.013695264444169
$row9[0,19] = -3.233880748928969
$row9[0,20] = 14.19083078680588
$row9[0,21] = -1.77235388068868
$row9[0,22] = 0.4853980060014154
$row9[0,23] = -11.61245515257498
$row9[0,24] = -20.00428048337838
$row9[0,25] = 3.852468463463529
$row9[0,26] = -10.91824497776777
$row9[0,27] = -27.42675730230798
$ws.Range("B9:AC9").Value = $row9

$row10 = New-Object 'object[,]' 1,28
$row10[0,0] = -1.2164722392424
$row10[0,1] = -11.45050024574915
$row10[0,2] = 5.379633434012145
$row10[0,3] = 1.057679041618949
$row10[0,4] = -5.537130777064892
$row10[0,5] = -8.400673372050916
$row10[0,6] = -19.05377386216486
$row10[0,7] = -14.55442616776373
$row10[0,8] = 11.95325502999418
$row10[0,9] = 1.864724338770341
$row10[0,10] = 1.431921363608129
$row10[0,11] = 8.1191116970252
$row10[0,12] = 11.4832281287722
$row10[0,13] = 2.975096678117097
$row10[0,14] = -4.716339417692327
$row10[0,15] = 4.318105369722222
$row10[0,16] = -5.074374610978587
$row10[0,17] = 1.670406697499622
$row10[0,18] = -3.117679282614945
$row10[0,19] = 1.512155377507596
$row10[0,20] = -3.822074686890403
$row10[0,21] = 1.728125437724763
$row10[0,22] = -2.578661058575849
$row10[0,23] = 7.416582675570321
$row10[0,24] = -2.173128685650406
$row10[0,25] = -7.109385930832257
$row10[0,26] = 0.4591684464186017
$row10[0,27] = -48.20252964268965
$ws.Range("B10:AC10").Value = $row10

$row11 = New-Object 'object[,]' 1,28
$row11[0,0] = 0.9288684790531327
$row11[0,1] = 7.41153211312891
$row11[0,2] = -12.53722374355275
$row11[0,3] = -13.12140189163336
$row11[0,4] = -5.378012006363437
$row11[0,5] = -8.705409258735472
$row11[0,6] = -6.828347155360569
$row11[0,7] = -0.7041355133930995
$row11[0,8] = 0.5244291281140994
$row11[0,9] = -5.85163348589084
$row11[0,10] = -1.619613858830224
$row11[0,11] = -16.12025393491743
$row11[0,12] = -20.73276359946449
$row11[0,13] = -2.05317976063329
$row11[0,14] = 3.333447273202541
$row11[0,15] = -11.06278207010978
$row11[0,16] = 13.90502491439268
$row11[0,17] = -1.551588885475272
$row11[0,18] = 4.939989175108254
$row11[0,19] = 1.081928499759411
$row11[0,20] = -5.652860501399569
$row11[0,21] = -3.641786494582059
$row11[0,22] = -2.974145992229329
$row11[0,23] = -19.52703952133388
$row11[0,24] = -13.82790031539912
$row11[0,25] = -25.81937297558348
$row11[0,26] = -22.27632385803323
$row11[0,27] = -35.46577332817661
$ws.Range("B11:AC11").Value = $row11
